$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Swap the "hommes" / "femmes" life-expectancy column headers (D2 <-> E2)
$d2 = $ws.Cells.Item(2, 4).Value2
$e2 = $ws.Cells.Item(2, 5).Value2
$ws.Cells.Item(2, 4).Value2 = $e2
$ws.Cells.Item(2, 5).Value2 = $d2

# 2) Swap the male/female life-expectancy data values between columns D and E
#    for every country data row (3-97)
for ($r = 3; $r -le 97; $r++) {
    $d = $ws.Cells.Item($r, 4).Value2
    $e = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 4).Value2 = $e
    $ws.Cells.Item($r, 5).Value2 = $d
}

# 3) Give the footer hyperlink rows (B109:B111) their own explicit cell style
#    (touching the font forces a distinct style record, matching the new
#    cellXfs entry introduced for these cells).
$footerLinks = $ws.Range("B109:B111")
$footerLinks.Font.Name = "Calibri"
